$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.546.17"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "2.638.93"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.52"
$ws.Range("E5").Value = "  -1.26%  "

$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("E7").Value = "  -1.16%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"
$ws.Range("E9").Value = "  -0.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.63"
$ws.Range("E10").Value = "  -3.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.96"
$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0815"
$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("E13").Value = "  +1.66%  "

$ws.Range("E14").Value = "  +3.60%  "

$ws.Range("D15").Value = "3.047.84"
$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("D16").Value = "2.633.53"
$ws.Range("E16").Value = "  -0.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.856"
$ws.Range("E17").Value = "  -1.80%  "

$ws.Range("D18").Value = "49.495.68"
$ws.Range("E18").Value = "  -0.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.39"
$ws.Range("E19").Value = "  +1.94%  "

$ws.Range("E20").Value = "  -0.41%  "

$ws.Range("E21").Value = "  -1.86%  "

$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "268.70"
$ws.Range("E23").Value = "  -3.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.10"
$ws.Range("E24").Value = "  -4.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("E25").Value = "  -0.65%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("E27").Value = "  -2.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("E28").Value = "  +1.46%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  -3.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.69"
$ws.Range("E31").Value = "  -3.83%  "

$ws.Range("E32").Value = "  -1.10%  "

$ws.Range("E33").Value = "  +0.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0815"
$ws.Range("E34").Value = "  +0.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.14"
$ws.Range("E35").Value = "  -1.96%  "

$ws.Range("E36").Value = "  -0.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.96"
$ws.Range("E37").Value = "  +2.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.03"
$ws.Range("E38").Value = "  -2.16%  "

$ws.Range("E39").Value = "  -0.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "129.52"
$ws.Range("E40").Value = "  +4.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.32"
$ws.Range("E41").Value = "  +3.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.02"
$ws.Range("E42").Value = "  +4.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0339"
$ws.Range("E43").Value = "  +7.53%  "

$ws.Range("E44").Value = "  -0.99%  "

$ws.Range("D45").Value = "2.062.74"
$ws.Range("E45").Value = "  -0.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.31"
$ws.Range("E46").Value = "  -0.36%  "

$ws.Range("E47").Value = "  +7.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.17"
$ws.Range("E48").Value = "  -7.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.88"
$ws.Range("E49").Value = "  -2.62%  "

$ws.Range("E50").Value = "  -2.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.64"
$ws.Range("E51").Value = "  -0.87%  "
